$wb = $excel.ActiveWorkbook

# Locate the most recent existing date-sheet (2025-07-22) so the new
# sheet is inserted immediately after it, at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "2025-07-23"

# Header row (rank / title / author / latest_episode), bold + thin border +
# centered-top alignment, matching the style used on every other date sheet.
$header = $ws.Range("A1:D1")
$headerValues = New-Object "object[,]" 1,4
$headerValues[0,0] = 'rank'
$headerValues[0,1] = 'title'
$headerValues[0,2] = 'author'
$headerValues[0,3] = 'latest_episode'
$header.Value = $headerValues
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

# Ranking data rows 2-51 (rank 1-50).
$data = New-Object "object[,]" 50,4
$data[0,0] = 1
$data[0,1] = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$data[0,2] = '光永康則'
$data[0,3] = '第６６話『六花停止』④'
$data[1,0] = 2
$data[1,1] = '地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。'
$data[1,2] = 'マツモトケンゴ'
$data[1,3] = '第６１話　偽彼氏の戦いが始まった'
$data[2,0] = 3
$data[2,1] = 'いとこのこ'
$data[2,2] = 'いぬちく(著者)'
$data[2,3] = '距離が近すぎる爽ちゃんと一緒に行きたい場所募集～！！【コメント募集企画】'
$data[3,0] = 4
$data[3,1] = '世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜'
$data[3,2] = '戸賀 環 坂木持丸 riritto'
$data[3,3] = '第49話②　城のパーティーに参加してみた'
$data[4,0] = 5
$data[4,1] = 'ラスボス討伐後に始める二周目冒険者ライフ はじまりの街でワケあり美少女たちがめちゃくちゃ懐いてきます'
$data[4,2] = '鬼麻正明(漫画) 朱月十話(原作) ファルまろ(キャラ原案)'
$data[4,3] = '第4話-2'
$data[5,0] = 6
$data[5,1] = '最果てのパラディン'
$data[5,2] = '奥橋睦（漫画） 柳野かなた（原作） 輪くすさが（キャラクター原案）'
$data[5,3] = '第66話　祝宴'
$data[6,0] = 7
$data[6,1] = '塔の管理をしてみよう'
$data[6,2] = '盧恩＆雪笠(Friendly Land)(著者) 早秋(原作) 雨神(キャラクター原案)'
$data[6,3] = '第91話前編'
$data[7,0] = 8
$data[7,1] = '絶対死なないステラ姫'
$data[7,2] = '光永康則 大高稲'
$data[7,3] = '第１４話　絶対旅立たない（２）'
$data[8,0] = 9
$data[8,1] = 'ひとりぼっちの異世界攻略'
$data[8,2] = 'びび（漫画） 五示正司（原作）'
$data[8,3] = '第228話　弱肉強食'
$data[9,0] = 10
$data[9,1] = 'Sランク冒険者である俺の娘たちは重度のファザコンでした'
$data[9,2] = 'しゅにち（漫画） 友橋かめつ（原作） 希望つばめ（原作イラスト）'
$data[9,3] = '第46話　敵か、味方か'
$data[10,0] = 11
$data[10,1] = '煽り系ゲーム配信者（20歳）、配信の切り忘れによりいい人バレする。'
$data[10,2] = '流嘉（漫画） 夏乃実（原作） 麦うさぎ（キャラクター原案）'
$data[10,3] = '第4話　サブ垢（後編）'
$data[11,0] = 12
$data[11,1] = '新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる'
$data[11,2] = '漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea'
$data[11,3] = '第69話'
$data[12,0] = 13
$data[12,1] = '江戸前エルフ'
$data[12,2] = '樋口彰彦'
$data[12,3] = '#115'
$data[13,0] = 14
$data[13,1] = 'ある日、惰眠を貪っていたら一族から追放されて森に捨てられました そのまま寝てたら周りが勝手に魔物の国を作ってたけど、私は気にせず今日も眠ります　コミック版'
$data[13,2] = '漫画/伊草さゆ 原作/白波ハクア キャラクター原案/まさよ'
$data[13,3] = 'chapter52【27話②】'
$data[14,0] = 15
$data[14,1] = '履いてください、鷹峰さん'
$data[14,2] = '柊裕一'
$data[14,3] = '第64話 当人同士で、どうぞ。(前編)'
$data[15,0] = 16
$data[15,1] = '聖女に嘘は通じない'
$data[15,2] = '日向 夏 浅見よう しんいし智歩'
$data[15,3] = '第25話②　あと二日で'
$data[16,0] = 17
$data[16,1] = 'ジャヒー様はくじけない！'
$data[16,2] = '昆布わかめ'
$data[16,3] = '復興計画その99 ジャヒー様と大切なあの日'
$data[17,0] = 18
$data[17,1] = '序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する'
$data[17,2] = '作画：マエD 原作：新人'
$data[17,3] = '第5話(1)'
$data[18,0] = 19
$data[18,1] = 'ダウナー系お姉さんに毎日カスの嘘を流し込まれる話'
$data[18,2] = '生倉のゑる(著者) はるばーど屋(原作者)'
$data[18,3] = '第11話'
$data[19,0] = 20
$data[19,1] = 'ゴリラ女子高生'
$data[19,2] = '大友しゅうま(著者)'
$data[19,3] = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'
$data[20,0] = 21
$data[20,1] = 'ラブコメと怪獣退治の不文律'
$data[20,2] = '御池慧（漫画） 上代務（原作） TMSLab（原作）'
$data[20,3] = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'
$data[21,0] = 22
$data[21,1] = '王子様の友達'
$data[21,2] = 'すけろく(著者)'
$data[21,3] = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'
$data[22,0] = 23
$data[22,1] = 'ゲーム　オブ　ファミリア-家族戦記-'
$data[22,2] = 'Ｄ．Ｐ(作画) 山口ミコト(原作)'
$data[22,3] = '第73話④'
$data[23,0] = 24
$data[23,1] = '傭兵団の料理番'
$data[23,2] = '梅木泰祐(漫画) 川井昂(原作) 四季童子(キャラクター原案)'
$data[23,3] = '第9話-2'
$data[24,0] = 25
$data[24,1] = 'みだりに憑かせてはなりません'
$data[24,2] = '栗田あぐり(著者)'
$data[24,3] = '第8話②'
$data[25,0] = 26
$data[25,1] = '宇崎ちゃんは遊びたい！'
$data[25,2] = '丈(著者)'
$data[25,3] = '第125話'
$data[26,0] = 27
$data[26,1] = '高森くんを黙らせたいっ!!'
$data[26,2] = '春乃カミハル'
$data[26,3] = '最終話'
$data[27,0] = 28
$data[27,1] = '勇者パーティを追放された【スキルサポーター】、仲間のスキルを解放して最強に成り上がる'
$data[27,2] = '作画：なかお 原作：前田氏'
$data[27,3] = '第6話(1)'
$data[28,0] = 29
$data[28,1] = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$data[28,2] = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$data[28,3] = '第５０話　雌雄を決する器用貧乏（３）'
$data[29,0] = 30
$data[29,1] = '米原くんはつよつよギャルから離れられない'
$data[29,2] = '川村拓(著者)'
$data[29,3] = '第15話'
$data[30,0] = 31
$data[30,1] = '悪役一家の奥方、死に戻りして心を入れ替える。'
$data[30,2] = '鏡(漫画) 丘野優(原作) TEDDY(キャラクター原案)'
$data[30,3] = '第32話①'
$data[31,0] = 32
$data[31,1] = 'ラーメン大好き小泉さん'
$data[31,2] = '鳴見なる'
$data[31,3] = '17杯目 家系'
$data[32,0] = 33
$data[32,1] = '最強の少年聖騎士、転生者を狩る'
$data[32,2] = '作画：御塩 原作：宇奈木ユラ'
$data[32,3] = '第6話(1)'
$data[33,0] = 34
$data[33,1] = '数分後の未来が分かるようになったけど、女心は分からない。'
$data[33,2] = 'You2(漫画) mty(原作)'
$data[33,3] = '第10話-2'
$data[34,0] = 35
$data[34,1] = 'オークの酒杯に祝福を'
$data[34,2] = 'かなどめはじめ'
$data[34,3] = '第44話　冥土黒子'
$data[35,0] = 36
$data[35,1] = '初歩魔法しか使わない謎の老魔法使いが旅をする'
$data[35,2] = '山代カゲツ(漫画) やまだのぼる(原作) にじまあるく(キャラクター原案)'
$data[35,3] = '第5話①'
$data[36,0] = 37
$data[36,1] = '8歳から始める魔法学'
$data[36,2] = 'ゆうなぎ（漫画） 上野夕陽（原作） 乃希（キャラクター原案）'
$data[36,3] = '第18話　意思交錯（前編）'
$data[37,0] = 38
$data[37,1] = '【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！'
$data[37,2] = '島知宏 音速炒飯 有都あらゆる'
$data[37,3] = '第２２食　ユクシーさんの覚悟、すごいのですわ！（３）'
$data[38,0] = 39
$data[38,1] = 'ちゃんと吸えない吸血鬼ちゃん'
$data[38,2] = '二式恭介(著者)'
$data[38,3] = '第101話：吸血鬼ちゃんのおまもり⑤'
$data[39,0] = 40
$data[39,1] = '宮廷魔導師、追放される　～無能だと追い出された最巧の魔導師は、部下を引き連れて冒険者クランを始めるようです～'
$data[39,2] = 'きつね丸（漫画） しんこせい（原作） ろこ（キャラクター原案）'
$data[39,3] = '第2話　憧れの人（後編）'
$data[40,0] = 41
$data[40,1] = 'ギルドを追放された回復術士、実は魔力無限だったので規格外の回復魔法で伝説となる'
$data[40,2] = '漫画：坂下コウ 原作：霞杏檎'
$data[40,3] = '第4話(1)'
$data[41,0] = 42
$data[41,1] = '二度追放された冒険者、激レアスキル駆使して美少女軍団を育成中！　コミック版'
$data[41,2] = '漫画/青木千尋 原作/南野雪花'
$data[41,3] = 'chapter67【34話②】'
$data[42,0] = 43
$data[42,1] = '転生したら没落貴族だったので、【呪言】を極めて家族を救います'
$data[42,2] = '作画：アマセケイ 原作：メソポ・たみあ'
$data[42,3] = '第6話(1)'
$data[43,0] = 44
$data[43,1] = '神猫ミーちゃんと猫用品召喚師の異世界奮闘記 ～目指すは、もふもふスローライフ！～'
$data[43,2] = 'にゃんたろう(原作) ねこのゆーま(作画) 岩崎美奈子(キャラクター原案)'
$data[43,3] = '第4話①'
$data[44,0] = 45
$data[44,1] = 'フシノカミ ～辺境から始める文明再生記～'
$data[44,2] = '黒杞よるの（漫画） 雨川水海（原作） 大熊まい（キャラクター原案）'
$data[44,3] = '第39話　蘇る歴史（前編）'
$data[45,0] = 46
$data[45,1] = 'ダンジョン・バスターズ　～中年男ですが庭にダンジョンが出現したので世界を救います～'
$data[45,2] = '蒼和 伸（漫画） 篠崎冬馬（原作） 千里GAN（キャラクター原案）'
$data[45,3] = '第32話　始動へ'
$data[46,0] = 47
$data[46,1] = '回復術士のやり直し'
$data[46,2] = '月夜涙(原作) 羽賀ソウケン(漫画) しおこんぶ(キャラクター原案)'
$data[46,3] = '第72話-1'
$data[47,0] = 48
$data[47,1] = '最凶貴族は死亡フラグを覆す'
$data[47,2] = '作画：sudekuma 原作：塚上'
$data[47,3] = '第6話(1)'
$data[48,0] = 49
$data[48,1] = 'ふかふかダンジョン攻略記～俺の異世界転生冒険譚～'
$data[48,2] = 'KAKERU'
$data[48,3] = '第66話「東アイギス」（後半）'
$data[49,0] = 50
$data[49,1] = '傷口と包帯'
$data[49,2] = '七井海星'
$data[49,3] = '第15話　新章開幕！　考える切谷'
$ws.Range("A2:D51").Value = $data

$ws.Range("A1").Select()
